# BOT; UPDATE DATA
# Adds the next day's (2020-05-11, serial 43962) case data to the three
# data sheets ("all", "kobe", "other"), pushing the trailing footnote
# row down by one, and updates the previous day's cumulative figures on
# "kobe" to reflect the newly confirmed case.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "kobe": bump cumulative totals for 2020-05-10 (row 88),
# then insert the new day's row (89) before the footnote row.
# ---------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")

$wsKobe.Range("D88").Value = 1
$wsKobe.Range("E88").Value = 278

$wsKobe.Rows("89:89").Insert()
$wsKobe.Range("A89").Value = 43962
$wsKobe.Range("B89").Value = 0
$wsKobe.Range("C89").Value = 2600
$wsKobe.Range("D89").Value = 0
$wsKobe.Range("E89").Value = 278
$wsKobe.Range("F89").Value = 75
$wsKobe.Range("G89").Value = 65
$wsKobe.Range("H89").Value = 10
$wsKobe.Range("I89").Value = 8
$wsKobe.Range("J89").Value = 177
$wsKobe.Range("A89").Select()

# ---------------------------------------------------------------
# Sheet "other": insert the new day's row (64) before the footnote row.
# ---------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")

$wsOther.Rows("64:64").Insert()
$wsOther.Range("A64").Value = 43962
$wsOther.Range("B64").Value = 0
$wsOther.Range("C64").Value = 14
$wsOther.Range("D64").Value = 5
$wsOther.Range("E64").Value = 4
$wsOther.Range("F64").Value = 1
$wsOther.Range("G64").Value = 0
$wsOther.Range("H64").Value = 9
$wsOther.Range("I64").Select()

# ---------------------------------------------------------------
# Sheet "all": insert the new day's row (34) before the footnote row.
# Handled last so that "all" ends up the active/selected sheet again,
# matching the original tabSelected state.
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

$wsAll.Rows("34:34").Insert()
$wsAll.Range("A34").Value = 43962
$wsAll.Range("B34").Value = 278
$wsAll.Range("C34").Value = 274
$wsAll.Range("D34").Value = 80
$wsAll.Range("E34").Value = 69
$wsAll.Range("F34").Value = 11
$wsAll.Range("G34").Value = 8
$wsAll.Range("H34").Value = 186
$wsAll.Range("A34").Select()
